$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 151, pushing existing rows 151-157 down to 152-158.
$ws.Rows.Item(151).Insert()

# Populate the newly inserted row 151 with the new weekly record.
$ws.Range("A151").Value = 4
$ws.Range("B151").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C151").Value = "Los Lagos"
$ws.Range("D151").Value = 44509
$ws.Range("E151").Value = 10
$ws.Range("F151").Value = 100112017
$ws.Range("G151").Value = "Apio"
$ws.Range("H151").Value = "Americana (o)"
$ws.Range("I151").Value = "Primera"
$ws.Range("J151").Value = 20
$ws.Range("K151").Value = 10000
$ws.Range("L151").Value = 10000
$ws.Range("M151").Value = 10000
$ws.Range("N151").Value = "$/docena de matas"
$ws.Range("O151").Value = "Región de Coquimbo"
$ws.Range("P151").Value = 1667
$ws.Range("Q151").Value = 6
$ws.Range("R151").Value = "Hortaliza"
